$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the as_of column keeps storing plain text (e.g. "Jan 2026") rather
# than being auto-converted into a date serial number by Excel.
$ws.Range("F8").NumberFormat = "@"
$ws.Range("F9").NumberFormat = "@"

# Row 8: Unemployment Rate (UNRATE) - update value, as_of date and recomputed metrics
$ws.Range("E8").Value = 4.3
$ws.Range("F8").Value = "Jan 2026"
$ws.Range("G8").Value = 4.581666666666664
$ws.Range("H8").Value = 0.2999999999999998
$ws.Range("I8").Value = 0.07499999999999996

# Row 9: Initial Jobless Claims (ICSA) - update value, as_of date and recomputed metrics
$ws.Range("E9").Value = 227000
$ws.Range("F9").Value = "Feb 2026"
$ws.Range("G9").Value = 363881.2260536398
$ws.Range("H9").Value = 5000
$ws.Range("I9").Value = 0.02252252252252252
